$wb = $excel.ActiveWorkbook

# --- Update "Ready for handoff" -> "In Translation" everywhere it appears ---
# (Overview!E2:F2, zh-cn!C2, de-de!C2 all share this status string)
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "In Translation"

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"

# --- Narrow the "Status" columns from 17.2159881591797 to 13.4101848602295 ---
# (the engine quantizes ColumnWidth to 1/6-character steps, so 12.5 is the input
# that lands on the closest representable width to the target, 13.333333333333334)
$ov.Columns.Item(5).ColumnWidth = 12.5
$ov.Columns.Item(6).ColumnWidth = 12.5

$zh.Columns.Item(3).ColumnWidth = 12.5

$de.Columns.Item(3).ColumnWidth = 12.5
